# Auto-generated edit script: update TPM-derived values for Lta-Tnfrsf1a sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.07057933333333334
$ws.Range("H2").Value = 0.211738
$ws.Range("I2").Value = 0.1281663575568867
$ws.Range("J2").Value = 0.1281663575568867
$ws.Range("M2").Value = 36.95112266666666
$ws.Range("N2").Value = 110.853368
$ws.Range("O2").Value = 0.1740115908809209
$ws.Range("P2").Value = 0.1775751473829744
$ws.Range("Q2").Value = 2.607985603731556
$ws.Range("R2").Value = 23.471870433584
$ws.Range("S2").Value = 0.0223024317758868
$ws.Range("T2").Value = 0.02275915983270315

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07057933333333334
$ws.Range("H3").Value = 0.211738
$ws.Range("I3").Value = 0.1281663575568867
$ws.Range("J3").Value = 0.1281663575568867
$ws.Range("O3").Value = 0.2442250025331967
$ws.Range("P3").Value = 0.2492264486514428
$ws.Range("Q3").Value = 3.660303819150445
$ws.Range("R3").Value = 32.942734372354
$ws.Range("S3").Value = 0.03130142899900125
$ws.Range("T3").Value = 0.03194244613049388

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.07057933333333334
$ws.Range("H4").Value = 0.211738
$ws.Range("I4").Value = 0.1281663575568867
$ws.Range("J4").Value = 0.1281663575568867
$ws.Range("M4").Value = 64.73785366666665
$ws.Range("N4").Value = 194.213561
$ws.Range("O4").Value = 0.3048658902295037
$ws.Range("P4").Value = 0.3111091917238571
$ws.Range("Q4").Value = 4.569154553224221
$ws.Range("R4").Value = 41.12239097901799
$ws.Range("S4").Value = 0.03907355069405314
$ws.Range("T4").Value = 0.0398737319057139

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.07057933333333334
$ws.Range("H5").Value = 0.211738
$ws.Range("I5").Value = 0.1281663575568867
$ws.Range("J5").Value = 0.1281663575568867
$ws.Range("M5").Value = 12.7841595
$ws.Range("N5").Value = 25.568319
$ws.Range("O5").Value = 0.06020363583370166
$ws.Range("P5").Value = 0.04095769119761797
$ws.Range("Q5").Value = 0.9022974547370002
$ws.Range("R5").Value = 5.413784728422001
$ws.Range("S5").Value = 0.007716080716486805
$ws.Range("T5").Value = 0.005249398094738457

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.07057933333333334
$ws.Range("H6").Value = 0.211738
$ws.Range("I6").Value = 0.1281663575568867
$ws.Range("J6").Value = 0.1281663575568867
$ws.Range("M6").Value = 46.01464833333333
$ws.Range("N6").Value = 138.043945
$ws.Range("O6").Value = 0.216693880522677
$ws.Range("P6").Value = 0.2211315210441077
$ws.Range("Q6").Value = 3.247683202934445
$ws.Range("R6").Value = 29.22914882641
$ws.Range("S6").Value = 0.02777286537145872
$ws.Range("T6").Value = 0.02834162159323733

# Row 7
$ws.Range("I7").Value = 0.05066414213561767
$ws.Range("J7").Value = 0.05066414213561768
$ws.Range("M7").Value = 36.95112266666666
$ws.Range("N7").Value = 110.853368
$ws.Range("O7").Value = 0.1740115908809209
$ws.Range("P7").Value = 0.1775751473829744
$ws.Range("Q7").Value = 1.0309363224
$ws.Range("R7").Value = 9.278426901599998
$ws.Range("S7").Value = 0.008816147973635929
$ws.Range("T7").Value = 0.008996692506764273

# Row 8
$ws.Range("I8").Value = 0.05066414213561767
$ws.Range("J8").Value = 0.05066414213561768
$ws.Range("O8").Value = 0.2442250025331967
$ws.Range("P8").Value = 0.2492264486514428
$ws.Range("S8").Value = 0.01237345024141346
$ws.Range("T8").Value = 0.01262684421843192

# Row 9
$ws.Range("I9").Value = 0.05066414213561767
$ws.Range("J9").Value = 0.05066414213561768
$ws.Range("M9").Value = 64.73785366666665
$ws.Range("N9").Value = 194.213561
$ws.Range("O9").Value = 0.3048658902295037
$ws.Range("P9").Value = 0.3111091917238571
$ws.Range("Q9").Value = 1.806186117299999
$ws.Range("R9").Value = 16.2556750557
$ws.Range("S9").Value = 0.01544576879488919
$ws.Range("T9").Value = 0.01576208030919463

# Row 10
$ws.Range("I10").Value = 0.05066414213561767
$ws.Range("J10").Value = 0.05066414213561768
$ws.Range("M10").Value = 12.7841595
$ws.Range("N10").Value = 25.568319
$ws.Range("O10").Value = 0.06020363583370166
$ws.Range("P10").Value = 0.04095769119761797
$ws.Range("Q10").Value = 0.35667805005
$ws.Range("R10").Value = 2.1400683003
$ws.Range("S10").Value = 0.003050165562959626
$ws.Range("T10").Value = 0.002075086288382854

# Row 11
$ws.Range("I11").Value = 0.05066414213561767
$ws.Range("J11").Value = 0.05066414213561768
$ws.Range("M11").Value = 46.01464833333333
$ws.Range("N11").Value = 138.043945
$ws.Range("O11").Value = 0.216693880522677
$ws.Range("P11").Value = 0.2211315210441077
$ws.Range("Q11").Value = 1.2838086885
$ws.Range("R11").Value = 11.5542781965
$ws.Range("S11").Value = 0.01097860956271946
$ws.Range("T11").Value = 0.01120343881284401

# Row 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.128589
$ws.Range("H12").Value = 0.385767
$ws.Range("I12").Value = 0.2335072176730087
$ws.Range("J12").Value = 0.2335072176730087
$ws.Range("M12").Value = 36.95112266666666
$ws.Range("N12").Value = 110.853368
$ws.Range("O12").Value = 0.1740115908809209
$ws.Range("P12").Value = 0.1775751473829744
$ws.Range("Q12").Value = 4.751507912584
$ws.Range("R12").Value = 42.763571213256
$ws.Range("S12").Value = 0.04063296242945773
$ws.Range("T12").Value = 0.04146507859327281

# Row 13
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.128589
$ws.Range("H13").Value = 0.385767
$ws.Range("I13").Value = 0.2335072176730087
$ws.Range("J13").Value = 0.2335072176730087
$ws.Range("O13").Value = 0.2442250025331967
$ws.Range("P13").Value = 0.2492264486514428
$ws.Range("Q13").Value = 6.668734111979001
$ws.Range("R13").Value = 60.01860700781101
$ws.Range("S13").Value = 0.05702830082771025
$ws.Range("T13").Value = 0.05819617459512336

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.128589
$ws.Range("H14").Value = 0.385767
$ws.Range("I14").Value = 0.2335072176730087
$ws.Range("J14").Value = 0.2335072176730087
$ws.Range("M14").Value = 64.73785366666665
$ws.Range("N14").Value = 194.213561
$ws.Range("O14").Value = 0.3048658902295037
$ws.Range("P14").Value = 0.3111091917238571
$ws.Range("Q14").Value = 8.324575865142998
$ws.Range("R14").Value = 74.92118278628699
$ws.Range("S14").Value = 0.07118838579089627
$ws.Range("T14").Value = 0.0726462417519365

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.128589
$ws.Range("H15").Value = 0.385767
$ws.Range("I15").Value = 0.2335072176730087
$ws.Range("J15").Value = 0.2335072176730087
$ws.Range("M15").Value = 12.7841595
$ws.Range("N15").Value = 25.568319
$ws.Range("O15").Value = 0.06020363583370166
$ws.Range("P15").Value = 0.04095769119761797
$ws.Range("Q15").Value = 1.6439022859455
$ws.Range("R15").Value = 9.863413715673001
$ws.Range("S15").Value = 0.01405798349732672
$ws.Range("T15").Value = 0.009563916513866053

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.128589
$ws.Range("H16").Value = 0.385767
$ws.Range("I16").Value = 0.2335072176730087
$ws.Range("J16").Value = 0.2335072176730087
$ws.Range("M16").Value = 46.01464833333333
$ws.Range("N16").Value = 138.043945
$ws.Range("O16").Value = 0.216693880522677
$ws.Range("P16").Value = 0.2211315210441077
$ws.Range("Q16").Value = 5.916977614535001
$ws.Range("R16").Value = 53.25279853081501
$ws.Range("S16").Value = 0.05059958512761768
$ws.Range("T16").Value = 0.05163580621880997

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.323617
$ws.Range("H17").Value = 0.970851
$ws.Range("I17").Value = 0.5876622826344869
$ws.Range("J17").Value = 0.5876622826344869
$ws.Range("M17").Value = 36.95112266666666
$ws.Range("N17").Value = 110.853368
$ws.Range("O17").Value = 0.1740115908809209
$ws.Range("P17").Value = 0.1775751473829744
$ws.Range("Q17").Value = 11.95801146401866
$ws.Range("R17").Value = 107.622103176168
$ws.Range("S17").Value = 0.1022600487019405
$ws.Range("T17").Value = 0.1043542164502342

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.323617
$ws.Range("H18").Value = 0.970851
$ws.Range("I18").Value = 0.5876622826344869
$ws.Range("J18").Value = 0.5876622826344869
$ws.Range("O18").Value = 0.2442250025331967
$ws.Range("P18").Value = 0.2492264486514428
$ws.Range("Q18").Value = 16.78305086062034
$ws.Range("R18").Value = 151.047457745583
$ws.Range("S18").Value = 0.1435218224650717
$ws.Range("T18").Value = 0.1464609837073936

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.323617
$ws.Range("H19").Value = 0.970851
$ws.Range("I19").Value = 0.5876622826344869
$ws.Range("J19").Value = 0.5876622826344869
$ws.Range("M19").Value = 64.73785366666665
$ws.Range("N19").Value = 194.213561
$ws.Range("O19").Value = 0.3048658902295037
$ws.Range("P19").Value = 0.3111091917238571
$ws.Range("Q19").Value = 20.95026999004566
$ws.Range("R19").Value = 188.552429910411
$ws.Range("S19").Value = 0.1791581849496651
$ws.Range("T19").Value = 0.1828271377570121

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.323617
$ws.Range("H20").Value = 0.970851
$ws.Range("I20").Value = 0.5876622826344869
$ws.Range("J20").Value = 0.5876622826344869
$ws.Range("M20").Value = 12.7841595
$ws.Range("N20").Value = 25.568319
$ws.Range("O20").Value = 0.06020363583370166
$ws.Range("P20").Value = 0.04095769119761797
$ws.Range("Q20").Value = 4.137171344911501
$ws.Range("R20").Value = 24.823028069469
$ws.Range("S20").Value = 0.03537940605692851
$ws.Range("T20").Value = 0.02406929030063061

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.323617
$ws.Range("H21").Value = 0.970851
$ws.Range("I21").Value = 0.5876622826344869
$ws.Range("J21").Value = 0.5876622826344869
$ws.Range("M21").Value = 46.01464833333333
$ws.Range("N21").Value = 138.043945
$ws.Range("O21").Value = 0.216693880522677
$ws.Range("P21").Value = 0.2211315210441077
$ws.Range("Q21").Value = 14.89112244968833
$ws.Range("R21").Value = 134.020102047195
$ws.Range("S21").Value = 0.1273428204608812
$ws.Range("T21").Value = 0.1299506544192164
